$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.551.41"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.591.97"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "659.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").Value = "3.589.62"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.93"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.203"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "4.259.97"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "95.442.17"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "3.578.53"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("E19").Value = "  -5.50%  "
$ws.Range("E20").Value = "  -7.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "511.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.13"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.39%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.89"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("D29").Value = "3.783.48"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.60"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.26"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.72"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +15.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.77"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +11.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.566"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "603.90"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.95%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +8.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.917"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "35.15"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +11.94%  "
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0420"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.43"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.35%  "
